# Add McKenzie basin thermal energy budget report, and switch to new 2010-19 actual
# weather data: update row 33 (Baseline 2010 12/16) and append new row 60
# (Baseline_2010-18 12/16) on the "2010 and 2010-18" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: label / year cells (no special formatting) ---
$ws.Range("B33").Value = 'Baseline 2010 12/16'
$ws.Range("C33").Value = 2010
$ws.Range("S33").Value = 2010

# --- Row 33: numeric cells (apply number format + highlight where changed) ---
$ws.Range("D33").Value = 831.51080300000001
$ws.Range("D33").NumberFormat = "0.00"
$ws.Range("D33").Interior.Color = 65535

$ws.Range("E33").Value = 1908.5467530000001
$ws.Range("E33").NumberFormat = "0.00"
$ws.Range("E33").Interior.Color = 65535

$ws.Range("F33").Value = 1.2760199999999999
$ws.Range("F33").NumberFormat = "0.00"

$ws.Range("G33").Value = 280.16485599999999
$ws.Range("G33").NumberFormat = "0.00"

$ws.Range("H33").Value = 10.610913999999999
$ws.Range("H33").NumberFormat = "0.00"

$ws.Range("I33").Value = 6.4750459999999999
$ws.Range("I33").NumberFormat = "0.00"
$ws.Range("I33").Interior.Color = 65535

$ws.Range("J33").Value = 8.8404570000000007
$ws.Range("J33").NumberFormat = "0.00"

$ws.Range("K33").Value = 737.20611599999995
$ws.Range("K33").NumberFormat = "0.00"
$ws.Range("K33").Interior.Color = 65535

$ws.Range("L33").Value = 59.834083999999997
$ws.Range("L33").NumberFormat = "0.00"
$ws.Range("L33").Interior.Color = 65535

$ws.Range("M33").Value = 1338.464966
$ws.Range("M33").NumberFormat = "0.00"
$ws.Range("M33").Interior.Color = 65535

$ws.Range("N33").Value = 895.30895999999996
$ws.Range("N33").NumberFormat = "0.00"
$ws.Range("N33").Interior.Color = 65535

$ws.Range("O33").Value = 6938.7304690000001
$ws.Range("O33").NumberFormat = "0"
$ws.Range("O33").Interior.Color = 65535

$ws.Range("P33").Value = 29450.638672000001
$ws.Range("P33").NumberFormat = "0"

$ws.Range("Q33").Value = 1.0701890000000001
$ws.Range("Q33").NumberFormat = "0.00"

$ws.Range("R33").Value = [double]"3.5199999999999999E-4"
$ws.Range("R33").NumberFormat = "0.000000"

# --- Row 60: label / year cells (new row, no special formatting) ---
$ws.Range("A60").Value = 'CW3M'
$ws.Range("B60").Value = 'Baseline_2010-18 12/16'
$ws.Range("C60").Value = '2010-18'
$ws.Range("S60").Value = '2010-18'

# --- Row 60: numeric cells (apply number format + highlight where changed) ---
$ws.Range("D60").Value = 897.84913466666671
$ws.Range("D60").NumberFormat = "0.00"
$ws.Range("D60").Interior.Color = 65535

$ws.Range("E60").Value = 1763.5263265555557
$ws.Range("E60").NumberFormat = "0.00"
$ws.Range("E60").Interior.Color = 65535

$ws.Range("F60").Value = 1.0174076666666665
$ws.Range("F60").NumberFormat = "0.00"

$ws.Range("G60").Value = 280.33542888888883
$ws.Range("G60").NumberFormat = "0.00"

$ws.Range("H60").Value = 9.775355222222224
$ws.Range("H60").NumberFormat = "0.00"

$ws.Range("I60").Value = 7.299440555555555
$ws.Range("I60").NumberFormat = "0.00"
$ws.Range("I60").Interior.Color = 65535

$ws.Range("J60").Value = 8.145128999999999
$ws.Range("J60").NumberFormat = "0.00"

$ws.Range("K60").Value = 646.63056122222224
$ws.Range("K60").NumberFormat = "0.00"

$ws.Range("L60").Value = 60.018756111111117
$ws.Range("L60").NumberFormat = "0.00"
$ws.Range("L60").Interior.Color = 65535

$ws.Range("M60").Value = 1342.5421007777777
$ws.Range("M60").NumberFormat = "0.00"
$ws.Range("M60").Interior.Color = 65535

$ws.Range("N60").Value = 902.73358833333339
$ws.Range("N60").NumberFormat = "0.00"
$ws.Range("N60").Interior.Color = 65535

$ws.Range("O60").Value = 5459.5160589999996
$ws.Range("O60").NumberFormat = "0"
$ws.Range("O60").Interior.Color = 65535

$ws.Range("P60").Value = 27227.338324888889
$ws.Range("P60").NumberFormat = "0"

$ws.Range("Q60").Value = 0.26704155555555559
$ws.Range("Q60").NumberFormat = "0.00"

$ws.Range("R60").Value = [double]"4.2666666666666656E-5"
$ws.Range("R60").NumberFormat = "0.000000"

# --- Update view: scroll / select to match the saved workbook view ---
$ws.Activate()
$ws.Range("O60").Select()
$excel.ActiveWindow.ScrollRow = 34

Write-Host "McKenzie regression rows updated"
